$wb = $excel.ActiveWorkbook

# This script applies a scheduled-runner data refresh to the per-job
# "Leve Profit" sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Each sheet holds
# market-price snapshots (columns H-N) for leve rows; values are plain
# numbers (no formulas) so cells are written directly. Where a computed
# profit cell (M/N) no longer applies it is cleared so it drops out of the
# sheet entirely, matching the source data feed.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1800
$ws.Range("I4").Value = 1750
$ws.Range("J4").Value = 1900
$ws.Range("K4").Value = 1750
$ws.Range("L4").Value = 1900
$ws.Range("M4").Value = -1636
$ws.Range("N4").Value = -2128
$ws.Range("H40").Value = 2750
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -2850
$ws.Range("H136").Value = 68609.86
$ws.Range("J136").Value = 68609.86
$ws.Range("L136").Value = 68609.86
$ws.Range("N136").Value = -78809.86
$ws.Range("H137").Value = 1707.625
$ws.Range("I137").Value = 1520
$ws.Range("J137").Value = 2420.6
$ws.Range("K137").Value = 4560
$ws.Range("L137").Value = 7261.799999999999
$ws.Range("M137").Value = -2010
$ws.Range("N137").Value = -12361.8
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = $null
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1450
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1450
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1450
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -1682
$ws.Range("H6").Value = 5005000
$ws.Range("I6").Value = 10000000
$ws.Range("K6").Value = 10000000
$ws.Range("M6").Value = -9999827
$ws.Range("H23").Value = 47431.145
$ws.Range("J23").Value = 38001.5
$ws.Range("L23").Value = 38001.5
$ws.Range("N23").Value = -38519.5
$ws.Range("H32").Value = 4214.634
$ws.Range("I32").Value = 2948.476
$ws.Range("K32").Value = 2948.476
$ws.Range("M32").Value = -2661.476
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = $null
$ws.Range("H74").Value = 1144.6
$ws.Range("I74").Value = 548.8570999999999
$ws.Range("K74").Value = 548.8570999999999
$ws.Range("M74").Value = 325.1429000000001
$ws.Range("H77").Value = 1144.6
$ws.Range("I77").Value = 548.8570999999999
$ws.Range("K77").Value = 2744.2855
$ws.Range("M77").Value = 1623.7145
$ws.Range("H132").Value = 2113.6
$ws.Range("I132").Value = 1669.6923
$ws.Range("K132").Value = 5009.0769
$ws.Range("M132").Value = -2479.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3817
$ws.Range("I20").Value = 3373.375
$ws.Range("K20").Value = 3373.375
$ws.Range("M20").Value = -3126.375
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = $null
$ws.Range("N38").Value = 0
$ws.Range("H82").Value = 25394.625
$ws.Range("I82").Value = 12539.25
$ws.Range("J82").Value = 38250
$ws.Range("K82").Value = 12539.25
$ws.Range("L82").Value = 38250
$ws.Range("M82").Value = -12156.25
$ws.Range("N82").Value = -39016
$ws.Range("H85").Value = 25394.625
$ws.Range("I85").Value = 12539.25
$ws.Range("J85").Value = 38250
$ws.Range("K85").Value = 12539.25
$ws.Range("L85").Value = 38250
$ws.Range("M85").Value = -11213.25
$ws.Range("N85").Value = -40902
$ws.Range("H105").Value = 2203.7036
$ws.Range("I105").Value = 2007.125
$ws.Range("K105").Value = 2007.125
$ws.Range("M105").Value = -260.125
$ws.Range("H107").Value = 782.45
$ws.Range("I107").Value = 554.8182
$ws.Range("K107").Value = 554.8182
$ws.Range("M107").Value = 1365.1818
$ws.Range("H108").Value = 20000
$ws.Range("J108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("N108").Value = -27680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("H22").Value = 1581.4546
$ws.Range("J22").Value = 1689.6
$ws.Range("L22").Value = 1689.6
$ws.Range("N22").Value = -2389.6
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250
$ws.Range("H69").Value = 12394
$ws.Range("I69").Value = 12394
$ws.Range("K69").Value = 12394
$ws.Range("M69").Value = -11645
$ws.Range("H72").Value = 12394
$ws.Range("I72").Value = 12394
$ws.Range("K72").Value = 37182
$ws.Range("M72").Value = -33438
$ws.Range("H107").Value = 891.3077
$ws.Range("I107").Value = 584
$ws.Range("K107").Value = 584
$ws.Range("M107").Value = 1336
$ws.Range("H132").Value = 2570.4614
$ws.Range("I132").Value = 1821.8
$ws.Range("J132").Value = 3038.375
$ws.Range("K132").Value = 5465.4
$ws.Range("L132").Value = 9115.125
$ws.Range("M132").Value = -2935.4
$ws.Range("N132").Value = -14175.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 30341.666
$ws.Range("J131").Value = 33025.91
$ws.Range("L131").Value = 99077.73000000001
$ws.Range("N131").Value = -109157.73
$ws.Range("H132").Value = 1206.8
$ws.Range("J132").Value = 1602.5
$ws.Range("L132").Value = 14422.5
$ws.Range("N132").Value = -19482.5
$ws.Range("H133").Value = 3784.2222
$ws.Range("J133").Value = 4999.8
$ws.Range("L133").Value = 14999.4
$ws.Range("N133").Value = -25119.4
$ws.Range("H134").Value = 2755.8965
$ws.Range("I134").Value = 1759.5
$ws.Range("J134").Value = 3982.2307
$ws.Range("K134").Value = 5278.5
$ws.Range("L134").Value = 11946.6921
$ws.Range("M134").Value = -208.5
$ws.Range("N134").Value = -22086.6921
$ws.Range("H136").Value = 1998.7778
$ws.Range("I136").Value = 1626.25
$ws.Range("J136").Value = 4979
$ws.Range("K136").Value = 4878.75
$ws.Range("L136").Value = 14937
$ws.Range("M136").Value = 221.25
$ws.Range("N136").Value = -25137
$ws.Range("H137").Value = 4677.6665
$ws.Range("J137").Value = 8999.666999999999
$ws.Range("L137").Value = 26999.001
$ws.Range("N137").Value = -37199.001
$ws.Range("H138").Value = 3116.5881
$ws.Range("J138").Value = 5076
$ws.Range("L138").Value = 15228
$ws.Range("N138").Value = -25508
$ws.Range("H140").Value = 1890.1
$ws.Range("I140").Value = 1038.4375
$ws.Range("J140").Value = 5296.75
$ws.Range("K140").Value = 3115.3125
$ws.Range("L140").Value = 15890.25
$ws.Range("M140").Value = 2064.6875
$ws.Range("N140").Value = -26250.25
$ws.Range("H141").Value = 5177.2
$ws.Range("I141").Value = 4713.25
$ws.Range("J141").Value = 7033
$ws.Range("K141").Value = 14139.75
$ws.Range("L141").Value = 21099
$ws.Range("M141").Value = -8959.75
$ws.Range("N141").Value = -31459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 9000
$ws.Range("J47").Value = 9000
$ws.Range("L47").Value = 9000
$ws.Range("N47").Value = -10136
$ws.Range("I70").Value = 4500
$ws.Range("K70").Value = 4500
$ws.Range("M70").Value = -4230
$ws.Range("I73").Value = 4500
$ws.Range("K73").Value = 4500
$ws.Range("M73").Value = -3564
$ws.Range("H97").Value = 1045.0834
$ws.Range("I97").Value = 1106.5
$ws.Range("K97").Value = 1106.5
$ws.Range("M97").Value = -610.5
$ws.Range("H102").Value = 2104.15
$ws.Range("I102").Value = 1949.2222
$ws.Range("J102").Value = 3498.5
$ws.Range("K102").Value = 1949.2222
$ws.Range("L102").Value = 3498.5
$ws.Range("M102").Value = -327.2221999999999
$ws.Range("N102").Value = -6742.5
$ws.Range("I132").Value = 19232668
$ws.Range("J132").Value = 4579.2
$ws.Range("K132").Value = 57698004
$ws.Range("L132").Value = 13737.6
$ws.Range("M132").Value = -57695474
$ws.Range("N132").Value = -18797.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1587.4667
$ws.Range("I46").Value = 996.2
$ws.Range("J46").Value = 1883.1
$ws.Range("K46").Value = 996.2
$ws.Range("L46").Value = 1883.1
$ws.Range("M46").Value = -808.2
$ws.Range("N46").Value = -2259.1
$ws.Range("H55").Value = 426.18518
$ws.Range("I55").Value = 344.8889
$ws.Range("K55").Value = 344.8889
$ws.Range("M55").Value = -171.8889
$ws.Range("H136").Value = 7547.375
$ws.Range("I136").Value = 5899
$ws.Range("K136").Value = 17697
$ws.Range("M136").Value = -15147

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9857
$ws.Range("J15").Value = 9857
$ws.Range("L15").Value = 9857
$ws.Range("N15").Value = -10433
$ws.Range("H113").Value = 444.5862
$ws.Range("I113").Value = 331
$ws.Range("K113").Value = 993
$ws.Range("M113").Value = 1177
